# Update "想去人数" (column F) values across the 展览 / 演出 / 全部类型 sheets
# to reflect the latest generated output (gh-pages data refresh).

$wb = $excel.ActiveWorkbook

# Sheet "展览": update column F ("想去人数") values
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 63
$ws.Cells.Item(3, 6).Value = 2872
$ws.Cells.Item(5, 6).Value = 6388
$ws.Cells.Item(6, 6).Value = 2485
$ws.Cells.Item(10, 6).Value = 2903
$ws.Cells.Item(11, 6).Value = 152
$ws.Cells.Item(12, 6).Value = 31
$ws.Cells.Item(13, 6).Value = 7154
$ws.Cells.Item(15, 6).Value = 23
$ws.Cells.Item(17, 6).Value = 223
$ws.Cells.Item(20, 6).Value = 8450
$ws.Cells.Item(24, 6).Value = 62
$ws.Cells.Item(26, 6).Value = 19
$ws.Cells.Item(28, 6).Value = 81
$ws.Cells.Item(31, 6).Value = 38
$ws.Cells.Item(32, 6).Value = 56
$ws.Cells.Item(33, 6).Value = 80
$ws.Cells.Item(34, 6).Value = 2602
$ws.Cells.Item(35, 6).Value = 42
$ws.Cells.Item(37, 6).Value = 32
$ws.Cells.Item(38, 6).Value = 1167
$ws.Cells.Item(40, 6).Value = 686
$ws.Cells.Item(41, 6).Value = 3683
$ws.Cells.Item(42, 6).Value = 6
$ws.Cells.Item(43, 6).Value = 180
$ws.Cells.Item(44, 6).Value = 19
$ws.Cells.Item(45, 6).Value = 1190
$ws.Cells.Item(46, 6).Value = 180
$ws.Cells.Item(47, 6).Value = 21
$ws.Cells.Item(48, 6).Value = 15

# Sheet "演出": update column F ("想去人数") values
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(2, 6).Value = 41
$ws.Cells.Item(4, 6).Value = 28
$ws.Cells.Item(5, 6).Value = 246
$ws.Cells.Item(7, 6).Value = 116

# Sheet "全部类型": update column F ("想去人数") values
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 41
$ws.Cells.Item(3, 6).Value = 2872
$ws.Cells.Item(4, 6).Value = 28
$ws.Cells.Item(5, 6).Value = 246
$ws.Cells.Item(6, 6).Value = 6388
$ws.Cells.Item(7, 6).Value = 2485
$ws.Cells.Item(8, 6).Value = 116
$ws.Cells.Item(12, 6).Value = 2903
$ws.Cells.Item(13, 6).Value = 153
$ws.Cells.Item(16, 6).Value = 31
$ws.Cells.Item(17, 6).Value = 7154
$ws.Cells.Item(20, 6).Value = 223
$ws.Cells.Item(23, 6).Value = 8450
$ws.Cells.Item(26, 6).Value = 62
$ws.Cells.Item(27, 6).Value = 19
$ws.Cells.Item(29, 6).Value = 81
$ws.Cells.Item(31, 6).Value = 38
$ws.Cells.Item(32, 6).Value = 56
$ws.Cells.Item(34, 6).Value = 80
$ws.Cells.Item(35, 6).Value = 2602
$ws.Cells.Item(36, 6).Value = 42
$ws.Cells.Item(38, 6).Value = 32
$ws.Cells.Item(39, 6).Value = 1167
$ws.Cells.Item(40, 6).Value = 686
$ws.Cells.Item(42, 6).Value = 3683
$ws.Cells.Item(43, 6).Value = 180
$ws.Cells.Item(44, 6).Value = 19
$ws.Cells.Item(46, 6).Value = 1190
$ws.Cells.Item(47, 6).Value = 180
$ws.Cells.Item(48, 6).Value = 21
$ws.Cells.Item(49, 6).Value = 15
